# Add a new "FAVÖK" (EBITDA) row of data to the DATA_RAW sheet (row 23)
# and to the "gelir tablosu (çeyreklik)" sheet (row 10), mirroring the
# upstream commit "Add files via upload".

$wb = $excel.ActiveWorkbook

# The same 64 period values (columns B..BM) are written to both sheets.
$values = @(
    0,0,0,-31362000,-37834000,0,-44246000,-54939000,-60777000,0,
    -84963000,-102475000,-101660000,0,11816000,8408000,-9354000,0,-15067000,7347000,
    1071000,0,-50346000,-109388000,-111862000,0,-149024000,-126140000,212422000,0,
    -97449000,-64316000,-337980000,0,-351774000,-369592000,716000000,0,-522000000,-167000000,
    261000000,0,283000000,327000000,809000000,0,594000000,1312000000,2004000000,0,
    3603000000,2464000000,3688000000,0,2664000000,2399000000,1764000000,0,1171000000,1605000000,
    2717000000,399000000,923000000,1624000000
)

$ws1 = $wb.Worksheets.Item("DATA_RAW")
$row1 = 23
$ws1.Cells.Item($row1, 1).Value = "FAVÖK"
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws1.Cells.Item($row1, $i + 2).Value = $values[$i]
}

$ws3 = $wb.Worksheets.Item("gelir tablosu (çeyreklik)")
$row3 = 10
$ws3.Cells.Item($row3, 1).Value = "FAVÖK"
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws3.Cells.Item($row3, $i + 2).Value = $values[$i]
}
